$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the "Count" column (currently column Y) so that
# "Count"/"Amount" shift one column to the right, making room for the new
# "DenominationName"/"DenominationValue" pair.
$ws.Columns.Item(25).Insert()

# Rename the old "ItemType" header (column X) to "DenominationName" and give
# the newly-inserted column (Y) the header "DenominationValue".
$ws.Range("X1").Value = "DenominationName"
$ws.Range("Y1").Value = "DenominationValue"
